# Weekly update: insert a new Albahaca price record for the
# "Terminal La Palmera de La Serena" market as row 34, pushing the
# existing historical rows (old 34-68) down by one (new 35-69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 34; all rows below shift down by one.
$ws.Rows("34:34").Insert()

# Populate the new row 34 with this week's record.
$ws.Range("A34").Value = 8
$ws.Range("B34").Value = "Terminal La Palmera de La Serena"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44539
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100112052
$ws.Range("G34").Value = "Albahaca"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 600
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = 3500
$ws.Range("N34").Value = "`$/paquete"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 3500
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"
